# Add a new "id" column (inserted before "version"), shift the "userAgent"
# column left to follow "os", and update the two data rows with new sample
# values: new profile names, newer user-agent strings, new GUID "id" values,
# a refreshed "created" timestamp, a "startURL" for the second profile, and
# the re-ordered "version"/"pathSave" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------------
$headers = @(
    "profileName", "browser", "os", "userAgent", "screen", "cpu",
    "languages", "startURL", "delayOpenSeconds", "webRTC", "getlocation",
    "timeZone", "clientRects", "audioContext", "fonts", "isRunning",
    "created", "id", "version", "pathSave", "proxy"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Row 2 -------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "Cuong KK9"
$ws.Cells.Item(2, 2).Value = "chrome"
$ws.Cells.Item(2, 3).Value = "'"
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 4).Value = "Mozilla/5.0 (Windows NT 10.0; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/132.0.0.0 Safari/537.36"
$ws.Cells.Item(2, 5).Value = "1920x1080"
$ws.Cells.Item(2, 6).Value = 4
$ws.Cells.Item(2, 7).Value = "vn"
$ws.Cells.Item(2, 8).Value = "'"
$ws.Cells.Item(2, 8).Style = "Normal"
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = "real"
$ws.Cells.Item(2, 11).Value = "prompt"
$ws.Cells.Item(2, 12).Value = "'"
$ws.Cells.Item(2, 12).Style = "Normal"
$ws.Cells.Item(2, 13).Value = "off"
$ws.Cells.Item(2, 14).Value = "off"
$ws.Cells.Item(2, 15).Value = "off"
$ws.Cells.Item(2, 16).Value = $false
$ws.Cells.Item(2, 17).Value = "15:10 7/2/25"
$ws.Cells.Item(2, 18).Value = "e5640726-544d-4243-9c2e-f51bf95df8a8"
$ws.Cells.Item(2, 19).Value = 130
$ws.Cells.Item(2, 20).Value = "E:\cuong-mmo\chromeProfile\Cuong KK9"

# --- Row 3 -------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "PHY 99"
$ws.Cells.Item(3, 2).Value = "chrome"
$ws.Cells.Item(3, 3).Value = "'"
$ws.Cells.Item(3, 3).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "Mozilla/5.0 (Windows NT 6.1; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/107.0.5304.107 Safari/537.36"
$ws.Cells.Item(3, 5).Value = "1920x1080"
$ws.Cells.Item(3, 6).Value = 4
$ws.Cells.Item(3, 7).Value = "vn"
$ws.Cells.Item(3, 8).Value = "http://zingnews.vn"
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = "real"
$ws.Cells.Item(3, 11).Value = "prompt"
$ws.Cells.Item(3, 12).Value = "'"
$ws.Cells.Item(3, 12).Style = "Normal"
$ws.Cells.Item(3, 13).Value = "off"
$ws.Cells.Item(3, 14).Value = "off"
$ws.Cells.Item(3, 15).Value = "off"
$ws.Cells.Item(3, 16).Value = $false
$ws.Cells.Item(3, 17).Value = "15:10 7/2/25"
$ws.Cells.Item(3, 18).Value = "dc0cdacb-e42b-4a60-8cfd-1d40dd65048a"
$ws.Cells.Item(3, 19).Value = 130
$ws.Cells.Item(3, 20).Value = "E:\cuong-mmo\chromeProfile\PHY 99"
